# Commit for AddAgency and Agency List
#
# Updates the "AddAgencyList" worksheet:
#   - Renames the shared "ATMNAgency..." test value
#   - Adds "Username" and "Password" columns (H1 / I1)
#   - Applies a wrap-text style to the existing G1 header cell
#   - Widens column G and increases the header row height
#   - Moves the active selection to J9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddAgencyList")

# Rename the existing agency test-data string in G2.
$ws.Range("G2").Value = "ATMNAgencybBx"

# Give the existing G1 header cell a wrap-text style (adds a new cellXfs entry).
$ws.Range("G1").WrapText = $true

# Add the two new header columns for the login flow.
$ws.Range("H1").Value = "Username"
$ws.Range("I1").Value = "Password"

# Make column G a bit wider to fit the header text.
$ws.Columns("G").ColumnWidth = 17.2

# Grow the header row to match the rest of the sheet's taller rows.
$ws.Rows("1").RowHeight = 30

# Leave the selection where the author left it when saving.
[void]$ws.Range("J9").Select()
